$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).Value = "'05.25.18"
    $ws.Cells.Item($r, 2).Value = "H.BROWN"
    $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
}
